$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Personalize the "I am here for you..." reply under handle_sadness with USER_NAME placeholder
$ws.Range("B2").Value = "I am here for you, USER_NAME. Together we can pass over it, ok?"

# Widen column B to fit the longer text (matches author's manual resize of the column)
$ws.Columns.Item(2).ColumnWidth = 61.25

# Move the active selection to B2 (the cell that was just edited)
$ws.Range("B2").Select() | Out-Null
